# Weighting & Scaling update & heatmap
# Updates the "Scaling" worksheet: adds two new columns (Optimal, Threshold)
# and refreshes the Min/Max scaling bounds for each mechanical parameter.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scaling")

# --- New header cells: E1 "Optimal", F1 "Threshold" ---
$ws.Range("E1").Value = "Optimal"
$ws.Range("F1").Value = "Threshold"

# Match the formatting used by the existing header cells (B1:D1)
$ws.Range("C1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Updated Min/Max values ---
# e_modulus
$ws.Range("B2").Value = 445
$ws.Range("C2").Value = 670

# tensile_strain_at_break
$ws.Range("B3").Value = 600
$ws.Range("C3").Value = 1929

# tensile_yield_strength
$ws.Range("B4").Value = 16
$ws.Range("C4").Value = 19.9

# --- Restore the UI selection state ---
[void]$ws.Range("F13").Select()
